# Append a new submissions row (row 4) to the sheet.
# All values -- including the numeric-looking ones -- must be stored as TEXT
# (matching the existing rows 2/3, which came in as inline/shared strings),
# so force the range to Text format before writing, then clear the
# formatting afterwards so no stray style index gets attached to the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = 4
$rowRange = $ws.Range("A$($targetRow):R$($targetRow)")

# Force text storage so "2", "0", "0.5" etc. are written as strings, not numbers.
$rowRange.NumberFormat = "@"

$values = @{
    "A" = "ggs"
    "B" = "hhshs"
    "C" = "hhsh"
    "D" = "hhshs"
    "E" = "2"
    "F" = "2"
    "G" = "2"
    "H" = "2"
    "I" = "2"
    "J" = "0"
    "K" = "0"
    "L" = "0"
    "M" = "0"
    "N" = "2"
    "O" = "0"
    "P" = "0"
    "Q" = "0.5"
    "R" = "2"
}

foreach ($col in $values.Keys) {
    $ws.Range("$col$targetRow").Value = $values[$col]
}

# Drop the temporary Text-format styling so the new cells end up unstyled,
# same as the rest of the data rows.
$rowRange.ClearFormats()
